$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F ("dSF") after a data repull / recalculation.
$updates = @{
    2  = -2
    4  = 1
    5  = 7
    6  = -2
    7  = -4
    8  = 9
    13 = 3
    14 = -2
    15 = -2
    16 = -2
    17 = -1
    18 = -2
    19 = 9
    20 = 2
    21 = -6
    22 = 11
    23 = 3
    24 = -2
    26 = -3
    27 = 1
    29 = 8
    30 = 11
    31 = 3
    32 = 0
    33 = 1
    36 = 4
    37 = 6
    38 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
